$d = $word.ActiveDocument

# Replace the three runs "Test update" + " #" + "5" with the new single
# sentence. Using wildcard find/replace across the whole paragraph text
# handles the run split transparently.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Test update #5",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "This is a test to see if the worddiff bot creates an .md file for diffing during PRs",
    2
)
